$wb = $excel.ActiveWorkbook

# --- Sheet: Save_Round ---
$wsRound = $wb.Worksheets.Item("Save_Round")
$wsRound.Range("D2").Value = 4
$wsRound.Range("E2").Value = 4
$wsRound.Range("F2").Value = 1
$wsRound.Range("G2").Value = "-"
$wsRound.Range("H2").Value = "-"
$wsRound.Range("J2").Value = "-"
$wsRound.Range("K2").Value = "-"
$wsRound.Range("L2").Value = 45562
$wsRound.Range("M2").Value = 0

# --- Sheet: Save_Holes ---
$wsHoles = $wb.Worksheets.Item("Save_Holes")
$wsHoles.Range("O2").Value = 1
$wsHoles.Rows.Item(3).Resize(2, 1).EntireRow.Delete()

# --- Sheet: Save_Shots ---
$wsShots = $wb.Worksheets.Item("Save_Shots")
$wsShots.Range("L2").Value = 1
$wsShots.Rows.Item(3).Resize(5, 1).EntireRow.Delete()
